$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the two runs in "var e = [] - 3; " into a single run with
#    a single trailing space (was: " e = [] - 3;" + " " as two runs).
# ------------------------------------------------------------------
$d.Content.Find.Execute(" e = [] - 3; ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " e = [] - 3; ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Append the new "TypeScript Syntax basics" section right before
#    the final (bookmarked) paragraph, preceded by a page break.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(1)  # wdCollapseStart

# Insert a paragraph mark + page break (its own paragraph).
$insertRange.InsertBreak(7)  # wdPageBreak

# Now insert all of the new textual paragraphs in one shot, separated
# by paragraph marks. The very last line merges into the (still) last
# paragraph of the document, right before the _GoBack bookmark - same
# as the other inserted paragraphs did.
$lines = @(
  "TypeScript Syntax basics:",
  "Dynamic vs Strong Typing",
  "JavaScript is dynamic",
  "var x;",
  "x = “foo”;",
  "I can then change the type of x",
  "x = 4;",
  "TypeScript is strongly typed",
  "var x : string;",
  "now we get a compile error on the second line",
  " Arrow functions",
  " add = (x: number, y: number) => { return x + y; }"
)
$text = [string]::Join("`r", $lines)

$insertRange = $lastPara.Range
$insertRange.Collapse(1)
$insertRange.InsertBefore($text)

# ------------------------------------------------------------------
# 3) Fix up the "var" prefix on the last line ("var add = ...") - it
#    needs a separate run boundary like the other "var" lines, so
#    re-insert it distinctly. Simplest: prefix text already lacks
#    "var" -> add it back through Find/Replace scoped to this phrase.
# ------------------------------------------------------------------
$d.Content.Find.Execute(" add = (x: number, y: number) => { return x + y; }", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "var add = (x: number, y: number) => { return x + y; }", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Apply "List Paragraph" style + bulleted list numbering to all of
#    the new bullet paragraphs (everything from "Dynamic vs Strong
#    Typing" through the final "var add = ..." paragraph).
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
# 12 lines were inserted; the last one merged into the pre-existing
# last paragraph, so 11 brand-new paragraphs were created, plus the
# "TypeScript Syntax basics:" paragraph before the list starts.
$typeScriptParaIndex = $n - 11
$firstListParaIndex = $typeScriptParaIndex + 1
$lastListParaIndex = $n

$startPara = $paras.Item($firstListParaIndex)
$endPara = $paras.Item($lastListParaIndex)
$listRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$listRange.Style = "List Paragraph"

$gallery = $word.ListGalleries.Item(1)
$template = $gallery.ListTemplates.Item(1)
$listRange.ListFormat.ApplyListTemplate($template)

# Bump the "level 1" (sub-bullet) items from ilvl 0 to ilvl 1.
$subIndexes = @(1, 2, 3, 4, 5, 6, 7, 8, 10)
foreach ($offset in $subIndexes) {
  $p = $paras.Item($firstListParaIndex + $offset)
  $p.Range.ListFormat.ListIndent()
}

# ------------------------------------------------------------------
# 5) Give the "List Paragraph" style the same look used by real Word
#    (uiPriority 34, left-indent of 720 twips / 36pt).
# ------------------------------------------------------------------
$style = $d.Styles("List Paragraph")
$style.Priority = 34
$style.ParagraphFormat.LeftIndent = 36
